$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Insert a new blank row above the existing row 17 (the old "Strings" header
# row, currently row 18) so a new Leetcode link can be slotted in right
# after the "Arrays" section (row 16, transpose-matrix) and before the
# "Strings" section header.
$ws.Rows.Item(17).EntireRow.Insert()

# Row insertion copies the formatting of the row above into every column of
# the new row (including column A), which would create a spurious styled
# blank cell that isn't part of the target layout. Only column C is used on
# this new row, so drop the inherited column-A cell entirely.
$ws.Range("A17").Clear()

# Keep the explicit row height consistent with the rest of the sheet.
$ws.Rows.Item(17).RowHeight = 15.75

# New Leetcode problem link for "Cells With Odd Values in a Matrix".
$newUrl = "https://leetcode.com/problems/cells-with-odd-values-in-a-matrix/"
$c17 = $ws.Range("C17")
$c17.Value = $newUrl
$ws.Hyperlinks.Add($c17, $newUrl) | Out-Null
# Hyperlinks.Add() stamps its own ad-hoc style onto the cell; reapply the
# sheet's shared "Hyperlink" cell style so it matches the other link cells.
$c17.Style = "Hyperlink"

# Restore the view: scroll position and the last active cell/selection.
$ws.Range("C26").Select()
